# Applies the LinkedIn-carousel text refresh described in the commit
# "Add generated LinkedIn draft": each of the 6 slides keeps its
# title / bullet / bullet layout, only the wording changes.
#
# NOTE: these text boxes use <a:bodyPr wrap="none"><a:spAutoFit/></a:bodyPr>,
# so writing TextRange.Text recalculates the shape's autofit Height as a
# side effect. We restore the original Height right after each shape's
# text is rewritten so only the wording changes, matching the source diff.

$p = $ppt.ActivePresentation

$slideText = @{
    1 = @(
        "Support for Manufacturing",
        "The renewable energy sector is seeking support for manufacturing.",
        "Focus on increasing domestic production capabilities."
    )
    2 = @(
        "Storage Solutions",
        "The sector is advocating for improvements in energy storage.",
        "Enhanced storage solutions are critical for renewable energy efficiency."
    )
    3 = @(
        "Green Hydrogen Initiatives",
        "There is a specific call for support in the development of green hydrogen.",
        "Green hydrogen is seen as a vital component for sustainable energy transition."
    )
    4 = @(
        "Overall Sector Goals",
        "The renewable energy sector aims to bolster its contributions to national energy goals.",
        "Support from the budget is seen as essential for future growth and innovation."
    )
    5 = @(
        "Importance of Policy Support",
        "Policy support is crucial for the sector to achieve its manufacturing goals.",
        "Investment in renewable technologies is expected to drive economic growth."
    )
    6 = @(
        "Collaboration with Government",
        "The renewable energy sector seeks collaboration with the government for strategic initiatives.",
        "Joint efforts are necessary to enhance the sector's infrastructure and capabilities."
    )
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $shape = $slide.Shapes.Item(1)
    $originalHeight = $shape.Height
    $originalWidth = $shape.Width

    $texts = $slideText[$i]
    $tr = $shape.TextFrame.TextRange
    for ($j = 1; $j -le $texts.Count; $j++) {
        $tr.Paragraphs($j).Runs(1).Text = $texts[$j - 1]
    }

    # Undo the autofit recalculation triggered by the text writes above.
    $shape.Width = $originalWidth
    $shape.Height = $originalHeight
}
